$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Column widths (C, D, F change; raw stored width = ColumnWidth + 5/6) ---
$ws.Columns.Item(3).ColumnWidth = (13 - 5/6)   # C: 16 -> 13
$ws.Columns.Item(4).ColumnWidth = (15 - 5/6)   # D: 13 -> 15
$ws.Columns.Item(6).ColumnWidth = (11 - 5/6)   # F: 15 -> 11

# --- Header row (month labels shift by one) ---
$ws.Range("C1").Value = "octubre"
$ws.Range("D1").Value = "noviembre"
$ws.Range("E1").Value = "diciembre"
$ws.Range("F1").Value = "enero"

# --- Data values ---
$ws.Range("C2").Value = 0

$ws.Range("D3").Value = 489.11
$ws.Range("E3").Value = 0

$ws.Range("C4").Value = 633.6
$ws.Range("D4").Value = 0

$ws.Range("C6").Value = 86.5
$ws.Range("D6").Value = 0

$ws.Range("C15").Value = 0

$ws.Range("C16").Value = 1382.33
$ws.Range("D16").Value = 0

$ws.Range("C17").Value = 0

$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 835.22
$ws.Range("E19").Value = 0

$ws.Range("C20").Value = 89.76000000000001
$ws.Range("D20").Value = 0

$ws.Range("D21").Value = 8791.290000000001
$ws.Range("E21").Value = 0

$ws.Range("C22").Value = 253.44
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 96.37

$ws.Range("C24").Value = 1391.04
$ws.Range("D24").Value = 4780.84
$ws.Range("E24").Value = 0

$ws.Range("C32").Value = 3836.67
$ws.Range("D32").Value = 14896.46
$ws.Range("E32").Value = 96.37
